# "Make trend use pollster analysis marked for it" / "make pollster
# analysis treat LIB/LNP in the same category."
#
# The sheet's row 2 (A2:I2) holds the raw first-preference figures typed
# in from a particular pollster's release; everything below (row 4 = %
# shares, row 6 = IFNA-cleaned % shares, K8 = swing) is formula-driven off
# those inputs, so re-pointing the analysis at a different pollster /
# re-bucketing LIB+LNP together just means re-typing this row of literals.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Avoid any intermediate/partial recalculation being baked into a cached
# formula result while we're mid-edit across several cells that all feed
# the same K2 total.
$excel.Calculation = -4135   # xlCalculationManual

# A=LNP B=ALP C=GRN D=ONP E=NXT F=UAP G=DEM H=DLP I=OTH
$ws.Range("A2").Value = 30
$ws.Range("B2").Value = 41
$ws.Range("C2").Value = "#N/A"
$ws.Range("D2").Value = "#N/A"
# E2 already holds the literal #N/A error and is left untouched.
$ws.Range("F2").Value = "#N/A"
$ws.Range("G2").Value = 16
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0

# Now let everything downstream (K2, row 4, row 6, K8) recompute off the
# refreshed inputs in one pass.
$excel.Calculation = -4105   # xlCalculationAutomatic
$excel.CalculateFullRebuild()

# The sheet was left with this range selected/active.
$ws.Range("A4:G4").Select()

$wb.Save()
